$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new fourth column "new" with its data
$ws.Range("D1").Value = "new"
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 0

# Move the active selection to E3 (matches the post-edit selection in the file)
$ws.Range("E3").Select() | Out-Null
